$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3838
$ws.Range("F5").Value = 1373
$ws.Range("F6").Value = 3834
$ws.Range("F7").Value = 391
$ws.Range("F8").Value = 201
$ws.Range("F10").Value = 8667
$ws.Range("F12").Value = 81
$ws.Range("F14").Value = 127
$ws.Range("F15").Value = 299
$ws.Range("F16").Value = 337
$ws.Range("F17").Value = 94
$ws.Range("F18").Value = 362
$ws.Range("F19").Value = 10987
$ws.Range("F20").Value = 292
$ws.Range("F23").Value = 189
$ws.Range("F24").Value = 12
$ws.Range("F25").Value = 159
$ws.Range("F29").Value = 2081
$ws.Range("F30").Value = 42
$ws.Range("F33").Value = 899
$ws.Range("F34").Value = 4087
$ws.Range("F36").Value = 284
$ws.Range("F38").Value = 3028
$ws.Range("F39").Value = 1248
$ws.Range("F40").Value = 174
$ws.Range("F42").Value = 347
$ws.Range("F43").Value = 331
$ws.Range("F44").Value = 46
$ws.Range("F45").Value = 111
$ws.Range("F46").Value = 128
$ws.Range("F47").Value = 90
$ws.Range("F48").Value = 99
$ws.Range("F49").Value = 85

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 36
$ws.Range("F19").Value = 175
$ws.Range("F22").Value = 54

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 12

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3838
$ws.Range("F6").Value = 1373
$ws.Range("F7").Value = 3834
$ws.Range("F8").Value = 391
$ws.Range("F10").Value = 201
$ws.Range("F11").Value = 8667
$ws.Range("F13").Value = 127
$ws.Range("F14").Value = 299
$ws.Range("F15").Value = 337
$ws.Range("F16").Value = 94
$ws.Range("F17").Value = 362
$ws.Range("F18").Value = 10987
$ws.Range("F19").Value = 292
$ws.Range("F23").Value = 189
$ws.Range("F24").Value = 12
$ws.Range("F25").Value = 36
$ws.Range("F26").Value = 159
$ws.Range("F30").Value = 2081
$ws.Range("F32").Value = 899
$ws.Range("F34").Value = 4087
$ws.Range("F36").Value = 284
$ws.Range("F38").Value = 3028
$ws.Range("F39").Value = 54
$ws.Range("F40").Value = 1248
$ws.Range("F41").Value = 174
$ws.Range("F43").Value = 347
$ws.Range("F44").Value = 331
$ws.Range("F45").Value = 111
$ws.Range("F46").Value = 128
$ws.Range("F47").Value = 90
$ws.Range("F48").Value = 99
$ws.Range("F49").Value = 85
